# Actualización automática 2025-06-02 14:06:09
# Adds a new "PRESUPUESTO" column (G) to the "VENTA MENSUAL" sheet,
# mirroring the formatting of the existing F column: bold/bordered header,
# currency-formatted zero values for each advisor row, and a right-aligned
# currency zero in the totals row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# New column width (stored OOXML width="17" <=> ColumnWidth 17 - 5/6)
$ws.Columns.Item(7).ColumnWidth = 17 - 5/6

# Header cell G1 - bold, centered/top aligned, thin border all around (matches F1)
$ws.Cells.Item(1, 7).Value = "PRESUPUESTO"
$ws.Cells.Item(1, 7).Font.Bold = $true
$ws.Cells.Item(1, 7).HorizontalAlignment = -4108
$ws.Cells.Item(1, 7).VerticalAlignment = -4160
$ws.Cells.Item(1, 7).Borders.LineStyle = 1

# Data rows G2:G6 - currency number format, value 0 (matches F2:F6)
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 7).Value = 0
    $ws.Cells.Item($r, 7).NumberFormat = '"$"#,##0.00'
}

# Totals row G7 - currency number format, right aligned, value 0 (matches F7)
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(7, 7).NumberFormat = '"$"#,##0.00'
$ws.Cells.Item(7, 7).HorizontalAlignment = -4152
